# AHP workbook update: add per-alternative Score column (J) and Rank column (K)
# to the Aggregation sheet, fixing the previous broken SUMPRODUCT/placeholder
# formulas, plus adding row-12 totals for the normalized columns F:I.

$wb = $excel.ActiveWorkbook

$agg = $wb.Worksheets.Item("Aggregation")

# --- Header row: insert a new "Score" column between Norm_Pekerjaan (I) and
#     Rank (now K instead of J) ---
$agg.Range("K1").Value = $agg.Range("J1").Value()
$agg.Range("J1").Value = "Score"

# --- Score (J) and Rank (K) formulas for each alternative row 2..11 ---
for ($r = 2; $r -le 11; $r++) {
    $agg.Range("J$r").Formula = "=(A17*F$r)+(B17*G$r)+(C17*H$r)+(D17*I$r)"
    $agg.Range("J$r").Style = "Normal"
    $agg.Range("K$r").Formula = "=RANK(J$r,J2:J11)"
    $agg.Range("K$r").Style = "Normal"
}

# --- Row 12 totals for the normalized weight columns ---
$agg.Range("F12").Formula = "=SUM(F2:F11)"
$agg.Range("G12").Formula = "=SUM(G2:G11)"
$agg.Range("H12").Formula = "=SUM(H2:H11)"
$agg.Range("I12").Formula = "=SUM(I2:I11)"

# --- Cosmetic: selection / active sheet state ---
$agg.Range("F32:F33").Select()
$agg.Activate()

$instr = $wb.Worksheets.Item("Instruksi_AHP")
$instr.Range("F6").Select()

$normalized = $wb.Worksheets.Item("Criteria_Normalized")
$normalized.Range("A1:H13").Select()

$altData = $wb.Worksheets.Item("Alternatives_Data")
$altData.Range("J9").Select()

$agg.Activate()
